$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 43: new date value
$ws.Range("D43").Value = 44783

# Update existing row 44: date + volume change
$ws.Range("D44").Value = 44757
$ws.Range("J44").Value = 300

# Append a new row 45 containing the data that used to live in row 44
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44391
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 100112043
$ws.Range("G45").Value = "Pepino dulce"
$ws.Range("H45").Value = "Cultivar IV Región"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = 15000
$ws.Range("N45").Value = "$/bandeja 18 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 833
$ws.Range("Q45").Value = 18
$ws.Range("R45").Value = "Hortaliza"

# Preserve the date-style formatting used by the rest of column D
$ws.Range("D45").NumberFormat = $ws.Range("D44").NumberFormat
